$wb = $excel.ActiveWorkbook

# Rename the "EActorType" sheet to "EActorClassType"
$wsEnum = $wb.Worksheets.Item("EActorType")
$wsEnum.Name = "EActorClassType"

$wsActor = $wb.Worksheets.Item("Actor")

# Update the Type field's enum reference
$wsActor.Range("C4").Value = "EActorClassType"

# Fill in "None" for the sample row's Type value
$wsActor.Range("C8").Value = "None"

# Add a new PrefabPath column (G)
$wsActor.Range("G2").Value = "All"
$wsActor.Range("G3").Value = "PrefabPath"
$wsActor.Range("G4").Value = "String"
$wsActor.Range("G8").Value = "Exported/Actor/001/Prefabs/Actor"

# Copy style from D4 (header style) onto G4
$wsActor.Range("D4").Copy()
$wsActor.Range("G4").PasteSpecial(-4122) | Out-Null
$wsActor.Range("G4").Value = "String"

# Column widths (widen to fit the longer "EActorClassType"/"PrefabPath" content).
# Column C grows to fit "EActorClassType" (~15.5 chars) and the new column G is
# sized to fit the exported prefab path (~32 chars). Excel quantizes ColumnWidth
# to whole pixels, so these inputs are chosen to land on the closest pixel step.
$wsActor.Columns.Item(3).ColumnWidth = 14.7
$wsActor.Columns.Item(7).ColumnWidth = 31.4

# Page setup for printing
$wsActor.PageSetup.PaperSize = 9
$wsActor.PageSetup.Orientation = 1

# Selections
$wsActor.Range("D10").Select() | Out-Null
$wsEnum.Range("B1").Select() | Out-Null

# Make Actor the active sheet/tab
$wsActor.Activate() | Out-Null
